$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 30,3
$data[0,0] = 10.46308040618896
$data[0,1] = -2.656880855560303
$data[0,2] = -2.153079509735107
$data[1,0] = 5.822948455810547
$data[1,1] = 0.93921160697937
$data[1,2] = -1.708410024642944
$data[2,0] = 76.90617370605469
$data[2,1] = -20.47752571105957
$data[2,2] = -10.06823921203613
$data[3,0] = 21.08987617492676
$data[3,1] = -20.16287040710449
$data[3,2] = 30.18997001647949
$data[4,0] = 10.72996425628662
$data[4,1] = -3.895856618881226
$data[4,2] = 4.558961868286133
$data[5,0] = 16.64472198486328
$data[5,1] = -21.42762756347656
$data[5,2] = 6.525417804718018
$data[6,0] = -22.31200981140137
$data[6,1] = 9.699769020080566
$data[6,2] = -17.87223434448242
$data[7,0] = -1.527544975280761
$data[7,1] = -2.107955932617188
$data[7,2] = -5.828549861907959
$data[8,0] = 25.63505172729492
$data[8,1] = -27.48458480834961
$data[8,2] = 41.41721725463867
$data[9,0] = 9.120732307434082
$data[9,1] = -3.210182666778564
$data[9,2] = 14.09486484527588
$data[10,0] = 10.57754039764404
$data[10,1] = 3.087275743484497
$data[10,2] = 21.89474105834961
$data[11,0] = -59.92874145507812
$data[11,1] = -16.7760124206543
$data[11,2] = -37.18990707397461
$data[12,0] = 18.34011650085449
$data[12,1] = 2.91422700881958
$data[12,2] = -21.99715805053711
$data[13,0] = 8.823348999023438
$data[13,1] = -6.199845790863037
$data[13,2] = -2.462417125701904
$data[14,0] = 36.45148849487305
$data[14,1] = -4.319328308105469
$data[14,2] = 22.25251770019531
$data[15,0] = 11.48440456390381
$data[15,1] = 5.108338356018066
$data[15,2] = 7.226218223571777
$data[16,0] = 1.909540176391602
$data[16,1] = -24.70354461669922
$data[16,2] = -9.863640785217283
$data[17,0] = 13.16304683685303
$data[17,1] = 1.566379547119141
$data[17,2] = -18.78874206542969
$data[18,0] = -26.10807037353516
$data[18,1] = -18.4659252166748
$data[18,2] = -26.81548118591309
$data[19,0] = 0.605715274810791
$data[19,1] = 0.4431395530700683
$data[19,2] = 13.34832572937012
$data[20,0] = -0.1363797187805175
$data[20,1] = 12.19958686828613
$data[20,2] = 1.961617946624756
$data[21,0] = 23.61079406738281
$data[21,1] = -17.57328033447266
$data[21,2] = -2.193590641021729
$data[22,0] = -16.4272632598877
$data[22,1] = 6.268374919891357
$data[22,2] = -29.44120216369629
$data[23,0] = -22.77685928344727
$data[23,1] = -11.73852729797363
$data[23,2] = -5.511547565460205
$data[24,0] = -20.23676681518555
$data[24,1] = -30.38065719604492
$data[24,2] = 13.72379970550537
$data[25,0] = -0.5074601173400879
$data[25,1] = 2.13600492477417
$data[25,2] = -3.543452739715576
$data[26,0] = 14.78357696533203
$data[26,1] = -19.61955070495605
$data[26,2] = 58.06728363037109
$data[27,0] = -78.68696594238281
$data[27,1] = -17.28630638122559
$data[27,2] = -26.55977630615234
$data[28,0] = -3.455702781677246
$data[28,1] = -1.554847240447998
$data[28,2] = -16.96383094787598
$data[29,0] = -10.43552684783936
$data[29,1] = -8.621312141418457
$data[29,2] = -5.604441642761231

$ws.Range("A2:C31").Value = $data
